$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Max_doners (column F) to 10 for rows 2-7
$ws.Range("F2:F7").Value = 10

# Update Run_experiment (column G) values per the target state
$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 1
$ws.Range("G6").Value = 0
$ws.Range("G7").Value = 0

# Move the active selection to G6
$ws.Range("G6").Select()
